$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.101987646160083
$ws.Range("D2").Value = 0.2360040478420729
$ws.Range("E2").Value = 0.1772031883258087
$ws.Range("F2").Value = 1.08873701812459
$ws.Range("G2").Value = 0.5299960472823386
$ws.Range("H2").Value = 0.6766214148207723
$ws.Range("I2").Value = 0.5084381643666029
$ws.Range("J2").Value = 0.1802673734551092
$ws.Range("K2").Value = 0.5884434309984954
$ws.Range("M2").Value = 0.245801753153394
$ws.Range("O2").Value = 2.382595943241114
$ws.Range("B3").Value = 0.0895997974957794
$ws.Range("D3").Value = 0.232933034930852
$ws.Range("E3").Value = 0.1767886261010148
$ws.Range("F3").Value = 1.093697044940349
$ws.Range("G3").Value = 0.533496198473685
$ws.Range("H3").Value = 0.6820436987439322
$ws.Range("I3").Value = 0.5164622036415363
$ws.Range("J3").Value = 0.1811333385756768
$ws.Range("K3").Value = 0.5137893058198699
$ws.Range("M3").Value = 0.2266639718715311
$ws.Range("O3").Value = 2.401113246222081
$ws.Range("B4").Value = 0.08197899534256692
$ws.Range("D4").Value = 0.2311365634280946
$ws.Range("E4").Value = 0.176615534805638
$ws.Range("F4").Value = 1.097359165632206
$ws.Range("G4").Value = 0.5360281235590065
$ws.Range("H4").Value = 0.6856782435459721
$ws.Range("I4").Value = 0.5216807868655025
$ws.Range("J4").Value = 0.1817584496640201
$ws.Range("K4").Value = 0.4677709483270007
$ws.Range("M4").Value = 0.2149394817383197
$ws.Range("O4").Value = 2.413924629891426
$ws.Range("B5").Value = 0.07887001799919346
$ws.Range("D5").Value = 0.2304270072126258
$ws.Range("E5").Value = 0.1765655197898646
$ws.Range("F5").Value = 1.099006636598041
$ws.Range("G5").Value = 0.5371560968597322
$ws.Range("H5").Value = 0.6872361539172545
$ws.Range("I5").Value = 0.5238807821590066
$ws.Range("J5").Value = 0.1820366874560904
$ws.Range("K5").Value = 0.4489740599635184
$ws.Range("M5").Value = 0.2101685842357952
$ws.Range("O5").Value = 2.419507844680652
$ws.Range("B6").Value = 0.0783535739577701
$ws.Range("D6").Value = 0.2303105492105288
$ws.Range("E6").Value = 0.1765584554168811
$ws.Range("F6").Value = 1.099289569151317
$ws.Range("G6").Value = 0.5373492037604848
$ws.Range("H6").Value = 0.6874994839428439
$ws.Range("I6").Value = 0.5242505219042983
$ws.Range("J6").Value = 0.1820843083356998
$ws.Range("K6").Value = 0.4458502307928143
$ws.Range("M6").Value = 0.2093768088477646
$ws.Range("O6").Value = 2.420456824143841
$ws.Range("B7").Value = 0.08193708019601331
$ws.Range("D7").Value = 0.2311269027802325
$ws.Range("E7").Value = 0.1766147771439641
$ws.Range("F7").Value = 1.097380755809837
$ws.Range("G7").Value = 0.5360429464112215
$ws.Range("H7").Value = 0.6856989430297773
$ws.Range("I7").Value = 0.5217101596470166
$ws.Range("J7").Value = 0.1817621069088879
$ws.Range("K7").Value = 0.4675176234366631
$ws.Range("M7").Value = 0.2148751112420157
$ws.Range("O7").Value = 2.413998459513053
$ws.Range("B8").Value = 0.09771948638559991
$ws.Range("D8").Value = 0.2349267224631149
$ws.Range("E8").Value = 0.177043365118152
$ws.Range("F8").Value = 1.090319289496271
$ws.Range("G8").Value = 0.531123403183706
$ws.Range("H8").Value = 0.6784276948164916
$ws.Range("I8").Value = 0.5111442969953099
$ws.Range("J8").Value = 0.1805465793184098
$ws.Range("K8").Value = 0.5627409124664666
$ws.Range("M8").Value = 0.2391978110071165
$ws.Range("O8").Value = 2.388681455346102
$ws.Range("B9").Value = 0.1285439155451229
$ws.Range("D9").Value = 0.2430814308972771
$ws.Range("E9").Value = 0.1785286099629602
$ws.Range("F9").Value = 1.081362473366376
$ws.Range("G9").Value = 0.5245172215503402
$ws.Range("H9").Value = 0.6665886767782041
$ws.Range("I9").Value = 0.4927393261901347
$ws.Range("J9").Value = 0.178903624972726
$ws.Range("K9").Value = 0.7479920436711325
$ws.Range("M9").Value = 0.2870896399018861
$ws.Range("O9").Value = 2.350477450978019
$ws.Range("B10").Value = 0.1511047983889142
$ws.Range("D10").Value = 0.2494968461965499
$ws.Range("E10").Value = 0.1800112336890436
$ws.Range("F10").Value = 1.077761828754475
$ws.Range("G10").Value = 0.521523032960971
$ws.Range("H10").Value = 0.659363091991807
$ws.Range("I10").Value = 0.4806271115228782
$ws.Range("J10").Value = 0.1781476758924256
$ws.Range("K10").Value = 0.8831387843072207
$ws.Range("M10").Value = 0.3223811823283356
$ws.Range("O10").Value = 2.329391053467049
$ws.Range("B11").Value = 0.1613477719853051
$ws.Range("D11").Value = 0.2525065639440385
$ws.Range("E11").Value = 0.1807703778828014
$ws.Range("F11").Value = 1.076770613755926
$ws.Range("G11").Value = 0.5205658519559364
$ws.Range("H11").Value = 0.65639518387691
$ws.Range("I11").Value = 0.4754228494086679
$ws.Range("J11").Value = 0.1779016570752816
$ws.Range("K11").Value = 0.9444020311814825
$ws.Range("M11").Value = 0.3384564620741699
$ws.Range("O11").Value = 2.32131588675125
$ws.Range("B12").Value = 0.16522342234002
$ws.Range("D12").Value = 0.253659301349856
$ws.Range("E12").Value = 0.1810699897956063
$ws.Range("F12").Value = 1.076488224038926
$ws.Range("G12").Value = 0.520261702497848
$ws.Range("H12").Value = 0.655317150311987
$ws.Range("I12").Value = 0.4734960767128662
$ws.Range("J12").Value = 0.1778225601745547
$ws.Range("K12").Value = 0.9675686362179761
$ws.Range("M12").Value = 0.3445464682463424
$ws.Range("O12").Value = 2.31847629660723
$ws.Range("B13").Value = 0.1643888754298786
$ws.Range("D13").Value = 0.253410461262078
$ws.Range("E13").Value = 0.1810049236650251
$ws.Range("F13").Value = 1.076544907829678
$ws.Range("G13").Value = 0.5203246115416533
$ws.Range("H13").Value = 0.6555472854620348
$ws.Range("I13").Value = 0.47390908590571
$ws.Range("J13").Value = 0.1778389696870093
$ws.Range("K13").Value = 0.9625807597520293
$ws.Range("M13").Value = 0.34323476463247
$ws.Range("O13").Value = 2.319078142032396
$ws.Range("B14").Value = 0.1616666884790021
$ws.Range("D14").Value = 0.2526011399712189
$ws.Range("E14").Value = 0.180794784039918
$ws.Range("F14").Value = 1.076745518510954
$ws.Range("G14").Value = 0.5205396601735117
$ws.Range("H14").Value = 0.6563055746979103
$ws.Range("I14").Value = 0.4752634513679697
$ws.Range("J14").Value = 0.1778948678862022
$ws.Range("K14").Value = 0.9463086204467004
$ws.Range("M14").Value = 0.3389574396236341
$ws.Range("O14").Value = 2.321077896122318
$ws.Range("B15").Value = 0.1599988537628718
$ws.Range("D15").Value = 0.2521071000312105
$ws.Range("E15").Value = 0.1806676473708428
$ws.Range("F15").Value = 1.076880503580441
$ws.Range("G15").Value = 0.5206789805959886
$ws.Range("H15").Value = 0.656776018740473
$ws.Range("I15").Value = 0.4760987658174747
$ws.Range("J15").Value = 0.1779309385986068
$ws.Range("K15").Value = 0.936337195716618
$ws.Range("M15").Value = 0.336337788239959
$ws.Range("O15").Value = 2.322331236407194
$ws.Range("B16").Value = 0.1504349701413474
$ws.Range("D16").Value = 0.2493019822771316
$ws.Range("E16").Value = 0.1799633227881579
$ws.Range("F16").Value = 1.077839617751359
$ws.Range("G16").Value = 0.5215937422673278
$ws.Range("H16").Value = 0.6595634692613999
$ws.Range("I16").Value = 0.4809733740998929
$ws.Range("J16").Value = 0.1781657224422446
$ws.Range("K16").Value = 0.879130623554488
$ws.Range("M16").Value = 0.3213310154929943
$ws.Range("O16").Value = 2.329949325146202
$ws.Range("B17").Value = 0.1445625125328007
$ws.Range("D17").Value = 0.2476044439739553
$ws.Range("E17").Value = 0.1795529069170421
$ws.Range("F17").Value = 1.078593621236138
$ws.Range("G17").Value = 0.5222586865331991
$ws.Range("H17").Value = 0.6613551742229618
$ws.Range("I17").Value = 0.4840420834388661
$ws.Range("J17").Value = 0.1783348149907411
$ws.Range("K17").Value = 0.843979975019181
$ws.Range("M17").Value = 0.3121299562728339
$ws.Range("O17").Value = 2.335011450496836
$ws.Range("B18").Value = 0.1411829547995609
$ws.Range("D18").Value = 0.2466366661318347
$ws.Range("E18").Value = 0.1793248203832398
$ws.Range("F18").Value = 1.079088180322692
$ws.Range("G18").Value = 0.5226792496963242
$ws.Range("H18").Value = 0.6624157474776169
$ws.Range("I18").Value = 0.4858358916104653
$ws.Range("J18").Value = 0.1784412846674712
$ws.Range("K18").Value = 0.8237420363846581
$ws.Range("M18").Value = 0.3068397472005486
$ws.Range("O18").Value = 2.338065843578704
$ws.Range("B19").Value = 0.1400383802430269
$ws.Range("D19").Value = 0.246310473756651
$ws.Range("E19").Value = 0.1792489648692488
$ws.Range("F19").Value = 1.079266086222134
$ws.Range("G19").Value = 0.5228281868885887
$ws.Range("H19").Value = 0.6627799982816924
$ws.Range("I19").Value = 0.4864481846065019
$ws.Range("J19").Value = 0.1784789159003353
$ws.Range("K19").Value = 0.8168863901321686
$ws.Range("M19").Value = 0.3050489275604633
$ws.Range("O19").Value = 2.339124529122003
$ws.Range("B20").Value = 0.1451878416084611
$ws.Range("D20").Value = 0.2477842604057798
$ws.Range("E20").Value = 0.1795957714019281
$ws.Range("F20").Value = 1.078507056088597
$ws.Range("G20").Value = 0.5221839577096858
$ws.Range("H20").Value = 0.6611613361914124
$ws.Range("I20").Value = 0.4837124366176022
$ws.Range("J20").Value = 0.1783158614382856
$ws.Range("K20").Value = 0.8477239242322696
$ws.Range("M20").Value = 0.3131092207544555
$ws.Range("O20").Value = 2.334457799123953
$ws.Range("B21").Value = 0.162466348045001
$ws.Range("D21").Value = 0.2528385046844761
$ws.Range("E21").Value = 0.1808561779561089
$ws.Range("F21").Value = 1.076684071640635
$ws.Range("G21").Value = 0.5204749118236407
$ws.Range("H21").Value = 0.6560816026808709
$ws.Range("I21").Value = 0.4748644479738582
$ws.Range("J21").Value = 0.1778780675718323
$ws.Range("K21").Value = 0.9510890345402743
$ws.Range("M21").Value = 0.3402137247817905
$ws.Range("O21").Value = 2.320484594184904
$ws.Range("B22").Value = 0.1737404232650874
$ws.Range("D22").Value = 0.2562176010562212
$ws.Range("E22").Value = 0.1817506727918961
$ws.Range("F22").Value = 1.076034480506991
$ws.Range("G22").Value = 0.5196978929992184
$ws.Range("H22").Value = 0.6530289314921305
$ws.Range("I22").Value = 0.4693380656440631
$ws.Range("J22").Value = 0.1776739222752042
$ws.Range("K22").Value = 1.018454260152339
$ws.Range("M22").Value = 0.3579433401748418
$ws.Range("O22").Value = 2.312624779731863
$ws.Range("B23").Value = 0.1677250066802571
$ws.Range("D23").Value = 0.2544072086683542
$ws.Range("E23").Value = 0.1812668033656308
$ws.Range("F23").Value = 1.07633161447508
$ws.Range("G23").Value = 0.5200814679955101
$ws.Range("H23").Value = 0.6546337579433725
$ws.Range("I23").Value = 0.4722641436590056
$ws.Range("J23").Value = 0.1777753801040944
$ws.Range("K23").Value = 0.9825180176620449
$ws.Range("M23").Value = 0.3484794370807833
$ws.Range("O23").Value = 2.316703233694682
$ws.Range("B24").Value = 0.1449051409160944
$ws.Range("D24").Value = 0.247702939973081
$ws.Range("E24").Value = 0.1795763678608431
$ws.Range("F24").Value = 1.078546001953335
$ws.Range("G24").Value = 0.5222176233853659
$ws.Range("H24").Value = 0.6612488753631851
$ws.Range("I24").Value = 0.4838613778510572
$ws.Range("J24").Value = 0.178324401506984
$ws.Range("K24").Value = 0.8460313761454756
$ws.Range("M24").Value = 0.3126664964930299
$ws.Range("O24").Value = 2.334707656028229
$ws.Range("B25").Value = 0.1202194636876754
$ws.Range("D25").Value = 0.2408005497037919
$ws.Range("E25").Value = 0.1780579548307379
$ws.Range("F25").Value = 1.083262105690906
$ws.Range("G25").Value = 0.5259782353922162
$ws.Range("H25").Value = 0.6695326579336154
$ws.Range("I25").Value = 0.4974706802011077
$ws.Range("J25").Value = 0.1792688341540796
$ws.Range("K25").Value = 0.6980409220084312
$ws.Range("M25").Value = 0.2741141460097083
$ws.Range("O25").Value = 2.359586973511469
